$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Weekly Timesheet" ---
$ws1 = $wb.Worksheets.Item("Weekly Timesheet")

# Row 2: Smithers/Regular day -> PTO day
$ws1.Range("B2").Value = "PTO"
$ws1.Range("C2").Value = 6.5
$ws1.Range("D2").Value = "PTO"
$ws1.Range("E2").Value = 150
$ws1.Range("F2").Value = 975

# Row 3: Bottomley -> Hewett
$ws1.Range("B3").Value = "Hewett"
$ws1.Range("C3").Value = 7
$ws1.Range("E3").Value = 150
$ws1.Range("F3").Value = 1050

# Row 4: Behrens -> Durfee
$ws1.Range("B4").Value = "Durfee"
$ws1.Range("C4").Value = 6
$ws1.Range("E4").Value = 150
$ws1.Range("F4").Value = 900

# Row 5: Goodrich -> Markfield
$ws1.Range("B5").Value = "Markfield"
$ws1.Range("C5").Value = 6
$ws1.Range("E5").Value = 150
$ws1.Range("F5").Value = 900

# Row 6: Campbell -> Corr
$ws1.Range("B6").Value = "Corr"
$ws1.Range("C6").Value = 6.5
$ws1.Range("E6").Value = 150
$ws1.Range("F6").Value = 975

# Row 8: SUBTOTAL row - hours drop 40 -> 32, total now populated
$ws1.Range("C8").Value = 32
$ws1.Range("D8").Value = "Reg: 32 / OT: 0"
$ws1.Range("F8").Value = 4800

# Row 12 / 13: ADMIN SUBTOTAL / GRAND TOTAL totals now populated
$ws1.Range("F12").Value = 4800
$ws1.Range("F13").Value = 4800

# --- Sheet 2: "Jason Schema" ---
$ws2 = $wb.Worksheets.Item("Jason Schema")

# Employee ID changed for every data row
$ws2.Range("B2:B6").Value = "emp_35u1tnme"

# Row 2: PTO day
$ws2.Range("D2").Value = "PTO"
$ws2.Range("E2").Value = 6.5
$ws2.Range("F2").Value = 150
$ws2.Range("G2").Value = 975
$ws2.Range("H2").Value = "PTO"
$ws2.Range("I2").Value = "PTO"

# Row 3: Bottomley -> Hewett
$ws2.Range("D3").Value = "Hewett"
$ws2.Range("E3").Value = 7
$ws2.Range("F3").Value = 150
$ws2.Range("G3").Value = 1050
$ws2.Range("H3").Value = "Regular"

# Row 4: Behrens -> Durfee
$ws2.Range("D4").Value = "Durfee"
$ws2.Range("E4").Value = 6
$ws2.Range("F4").Value = 150
$ws2.Range("G4").Value = 900
$ws2.Range("H4").Value = "Regular"

# Row 5: Goodrich -> Markfield
$ws2.Range("D5").Value = "Markfield"
$ws2.Range("E5").Value = 6
$ws2.Range("F5").Value = 150
$ws2.Range("G5").Value = 900
$ws2.Range("H5").Value = "Regular"

# Row 6: Campbell -> Corr
$ws2.Range("D6").Value = "Corr"
$ws2.Range("E6").Value = 6.5
$ws2.Range("F6").Value = 150
$ws2.Range("G6").Value = 975
$ws2.Range("H6").Value = "Regular"
